$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.829.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.45%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.570.36'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.23%  '

$ws.Range('E4').Value = '  +0.31%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '621.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.07%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.45'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.25%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.570.77'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.18%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.492'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.148'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.36'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.45%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.439'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000224'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.36%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.34'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.09%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.181.50'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.45%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.341.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.37%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.582.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.95%  '

$ws.Range('E18').Value = '  +0.53%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.79'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.41%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.04'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.84%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.07'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.73%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '459.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.52%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.642'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.27%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000132'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.28%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.61'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.62%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.726.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.63%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.17%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.89%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.61'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.20%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.98%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.172'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.57%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.35%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.46'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.49%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.55%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.92'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.78%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.570.71'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.76%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.33'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.24%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.84%  '

$ws.Range('E40').Value = '  -0.05%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '179.82'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.71%  '

$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.52%  '

$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0919'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.73%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.69%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '30.90'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +17.94%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.904'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.38%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.37'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.39%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '45.75'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.22%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.71'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.47%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.11%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.264'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.87%  '
